# Slide 10 ("Auswirkungen auf Marketing"): rename the product "Emotion
# Detector" to "ShireEye" inside the "Leitfragen:" bullet list, and
# re-split every paragraph's runs the way PowerPoint does after an
# in-place text edit (one run per word / whitespace token, formatting
# unchanged throughout).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Re-applying the paragraph's own (unchanged) font size to a sub-range
# forces PowerPoint's run-splitting logic without altering how the text
# looks.
function Split-Runs($para, $ranges) {
    $fontSize = $para.Font.Size
    foreach ($r in $ranges) {
        $chunk = $para.Characters($r[0], $r[1])
        $chunk.Font.Size = $fontSize
    }
}

# --- Paragraph 1: "Leitfragen:" -> "Leitfragen" + ":" ----------------------
$para1 = $tr.Paragraphs(1)
Split-Runs $para1 @(, @(1, 10), @(11, 1))

# --- Paragraph 2 --------------------------------------------------------
# "Wie lässt sich der Emotion Detector erfolgreich als Kaufargument vermarkten?"
$para2 = $tr.Paragraphs(2)
Split-Runs $para2 @(
    , @(1, 4), @(5, 5), @(10, 1), @(11, 4), @(15, 1), @(16, 20), @(36, 1),
    @(37, 11), @(48, 1), @(49, 3), @(52, 1), @(53, 12), @(65, 1), @(66, 10), @(76, 1)
)
# Replace the now-isolated "der Emotion Detector" run with "ShireEye".
$para2.Characters(16, 20).Text = "ShireEye"

# --- Paragraph 4 --------------------------------------------------------
# "Was für Auswirkungen hat die Einführung des Emotion Detectors auf den Absatz der Bankprodukte?"
$para4 = $tr.Paragraphs(4)
Split-Runs $para4 @(
    , @(1, 4), @(5, 3), @(8, 1), @(9, 12), @(21, 9), @(30, 10), @(40, 1),
    @(41, 21), @(62, 9), @(71, 6), @(77, 5), @(82, 12), @(94, 1)
)
# Replace right-to-left so earlier character offsets stay valid.
$para4.Characters(41, 21).Text = "ShireEye"
$para4.Characters(40, 1).Text = " von "

# --- Paragraph 6 (text unchanged, only re-split into runs) --------------
# "Welche Zweifel hegen Zielgruppen und wie räumt man diese aus? "
$para6 = $tr.Paragraphs(6)
Split-Runs $para6 @(
    , @(1, 6), @(7, 1), @(8, 7), @(15, 1), @(16, 5), @(21, 1), @(22, 11),
    @(33, 5), @(38, 3), @(41, 1), @(42, 5), @(47, 5), @(52, 5), @(57, 1), @(58, 3), @(61, 2)
)
